$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "similarity_score" column (F),
# pushing it to column H. The inserted columns F and G will hold the new
# BLEU_score and ROUGE_score metrics.
$ws.Columns("F:G").Insert()

# Header row
$ws.Range("F1").Value = "BLEU_score"
$ws.Range("G1").Value = "ROUGE_score"

# BLEU_score values (column F)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0.1207742128058413
$ws.Range("F4").Value = [double]"6.373704167435469E-155"
$ws.Range("F5").Value = [double]"6.373704167435469E-155"
$ws.Range("F6").Value = 0.3719426400061713

# ROUGE_score values (column G)
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0.2916666666666666
$ws.Range("G4").Value = 0.2222222222222222
$ws.Range("G5").Value = 0.2222222222222222
$ws.Range("G6").Value = 0.6896551724137931
